# Install Bare Equipment Template Complete
#
# 1) "${TITLE}" -> "${TYPE}" in the Recommendation heading, split
#    across three runs: "${T" / "YPE" / "}"
# 2) " ${TEMPS" + "TR" -> " ${TEMP" + "S" in the "Current Practice and
#    Observations" paragraph (the run holding the trailing "}" is left
#    untouched).
#
# Both edits are applied via Range.InsertXML: calling InsertXML on a
# Range that still spans the *original* text (i.e. not pre-collapsed)
# replaces that span in place and keeps the surrounding runs/paragraph
# untouched, which is exactly the kind of surgical run-split the diff
# calls for.

$d = $word.ActiveDocument

function New-PkgXml([string]$bodyXml) {
    return '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
           $bodyXml +
           '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------
# Change 1: ${TITLE} -> ${T}/YPE/} split runs
# ---------------------------------------------------------------------
$rngTitle = $d.Content
$foundTitle = $rngTitle.Find.Execute("`${TITLE}", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundTitle) {
    throw "Could not find `${TITLE} placeholder"
}

$titleBody = '<w:p><w:r><w:t>${T</w:t></w:r><w:r><w:t>YPE</w:t></w:r><w:r><w:t>}</w:t></w:r></w:p>'
$titleRange = $d.Range($rngTitle.Start, $rngTitle.End)
$titleRange.InsertXML((New-PkgXml $titleBody))

# ---------------------------------------------------------------------
# Change 2: ${TEMPSTR} -> ${TEMPS} (mid-paragraph, keeps trailing "}"
#           run untouched)
# ---------------------------------------------------------------------
$rngTemp = $d.Content
$foundTemp = $rngTemp.Find.Execute(" `${TEMPSTR}", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundTemp) {
    throw "Could not find `${TEMPSTR} placeholder"
}

$tempPara = $rngTemp.Paragraphs(1)
$paraEnd = $tempPara.Range.End

$tempBody = '<w:p w14:paraId="6E433836" w14:textId="2488BB27" w:rsidR="00D935C7" w:rsidRPr="00471BE2" w:rsidRDefault="00D935C7" w:rsidP="00D935C7">' +
            '<w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:b/></w:rPr></w:pPr>' +
            '<w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> ${TEMP</w:t></w:r>' +
            '<w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>S</w:t></w:r>' +
            '<w:r w:rsidR="007559E7"><w:rPr><w:color w:val="000000"/></w:rPr><w:t>}</w:t></w:r>' +
            '<w:r w:rsidRPr="009725DB"><w:rPr><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve">. The exposed surfaces on these are hot enough to cause injury when the </w:t></w:r>' +
            '<w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve">${TYPE} </w:t></w:r>' +
            '<w:r w:rsidRPr="009725DB"><w:rPr><w:color w:val="000000"/></w:rPr><w:t>are working. Insulating these will also reduce heat losses through the convection.</w:t></w:r>' +
            '</w:p>'

$tempRange = $d.Range($rngTemp.Start, $paraEnd)
$tempRange.InsertXML((New-PkgXml $tempBody))

Write-Output "done"
